$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = 2.2
$ws.Range("R2").Value = 1.67
$ws.Range("BD2").Value = 126
$ws.Range("G3").Value = 2.35
$ws.Range("I3").Value = 3.3
$ws.Range("J3").Value = 3.1
$ws.Range("K3").Value = 2.05
$ws.Range("L3").Value = 3.75
$ws.Range("U3").Value = 1.91
$ws.Range("V3").Value = 1.91
$ws.Range("W3").Value = 7
$ws.Range("X3").Value = 11
$ws.Range("AE3").Value = 15
$ws.Range("AG3").Value = 301
$ws.Range("AH3").Value = 9
$ws.Range("AV3").Value = 51
$ws.Range("BA3").Value = 81
$ws.Range("BD3").Value = 151
$ws.Range("G6").Value = 1.66
$ws.Range("M6").Value = 1.05
$ws.Range("N6").Value = 11
$ws.Range("Q6").Value = 2
$ws.Range("R6").Value = 1.85
$ws.Range("G7").Value = 3.3
$ws.Range("H7").Value = 3.2
$ws.Range("I7").Value = 2.15
$ws.Range("J7").Value = 3.75
$ws.Range("K7").Value = 2.1
$ws.Range("M7").Value = 1.06
$ws.Range("N7").Value = 10
$ws.Range("O7").Value = 1.3
$ws.Range("P7").Value = 3.4
$ws.Range("Q7").Value = 2.03
$ws.Range("R7").Value = 1.83
$ws.Range("U7").Value = 1.8
$ws.Range("V7").Value = 1.95
$ws.Range("W7").Value = 10
$ws.Range("AC7").Value = 9.5
$ws.Range("AE7").Value = 15
$ws.Range("AF7").Value = 51
$ws.Range("AG7").Value = 251
$ws.Range("AH7").Value = 7.5
$ws.Range("AL7").Value = 19
$ws.Range("AM7").Value = 29
$ws.Range("AP7").Value = 26
$ws.Range("AR7").Value = 81
$ws.Range("AS7").Value = 201
$ws.Range("AU7").Value = 8
$ws.Range("AX7").Value = 13
$ws.Range("AY7").Value = 23
$ws.Range("BA7").Value = 67
$ws.Range("G8").Value = 1.86
$ws.Range("Q8").Value = 1.7
$ws.Range("R8").Value = 2.1
$ws.Range("Q9").Value = 2.25
$ws.Range("R9").Value = 1.62
$ws.Range("M10").Value = 1.07
$ws.Range("N10").Value = 9
$ws.Range("G11").Value = 2.25
$ws.Range("H11").Value = 3.2
$ws.Range("I11").Value = 3.2
$ws.Range("J11").Value = 3
$ws.Range("K11").Value = 2.1
$ws.Range("M11").Value = 1.07
$ws.Range("N11").Value = 9
$ws.Range("U11").Value = 1.8
$ws.Range("V11").Value = 1.91
$ws.Range("AA11").Value = 19
$ws.Range("AC11").Value = 9
$ws.Range("AH11").Value = 9.5
$ws.Range("AS11").Value = 151
$ws.Range("AY11").Value = 26
$ws.Range("G13").Value = 2.4
$ws.Range("I13").Value = 3.2
$ws.Range("Q13").Value = 2.35
$ws.Range("R13").Value = 1.57
$ws.Range("U13").Value = 2
$ws.Range("V13").Value = 1.73
$ws.Range("AC13").Value = 7
$ws.Range("AE13").Value = 17
$ws.Range("AS13").Value = 251
$ws.Range("AZ13").Value = 67
$ws.Range("BA13").Value = 101
$ws.Range("Q20").Value = 2.03
$ws.Range("R20").Value = 1.83
$ws.Range("Q21").Value = 2.35
$ws.Range("R21").Value = 1.57
$ws.Range("Q23").Value = 2
$ws.Range("R23").Value = 1.85
$ws.Range("I25").Value = 7.1
$ws.Range("J25").Value = 1.85
$ws.Range("L25").Value = 6.2
$ws.Range("P25").Value = 4.15
$ws.Range("R25").Value = 2.22
$ws.Range("S25").Value = 1.3
$ws.Range("U25").Value = 1.8
$ws.Range("V25").Value = 1.9
$ws.Range("W25").Value = 7.1
$ws.Range("X25").Value = 7.6
$ws.Range("Y25").Value = 8.5
$ws.Range("AA25").Value = 11.75
$ws.Range("AB25").Value = 26
$ws.Range("AE25").Value = 18.5
$ws.Range("AF25").Value = 80
$ws.Range("AG25").Value = 600
$ws.Range("AH25").Value = 18.5
$ws.Range("AI25").Value = 55
$ws.Range("AN25").Value = 3.4
$ws.Range("AP25").Value = 14.5
$ws.Range("AQ25").Value = 17
$ws.Range("AU25").Value = 7.6
